$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71; everything from row 71 downward (old rows
# 71-170) shifts down to 72-171, matching the target diff (dimension becomes
# A1:T171 and old row 170's data lands on new row 171).
$ws.Rows.Item(71).Insert()

# Populate the freshly inserted row 71 with the new record.
$ws.Range("A71").Value = 3
$ws.Range("B71").Value = "Femacal de La Calera"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44546
$ws.Range("E71").Value = 5
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100101
$ws.Range("H71").Value = "Berries"
$ws.Range("I71").Value = 100101001
$ws.Range("J71").Value = "Arándano (blue)"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 40
$ws.Range("N71").Value = 6000
$ws.Range("O71").Value = 6000
$ws.Range("P71").Value = 6000
$ws.Range("Q71").Value = "$/bandeja 2 kilos"
$ws.Range("R71").Value = "Provincia de Cardenal Caro"
$ws.Range("S71").Value = 3000
$ws.Range("T71").Value = 2
